# Saldo.xlsx data refresh:
#  - 3 existing accounts receive updated Saldo (balance) values
#  - 1 brand-new account/row is added
#  - The whole table is re-sorted by Saldo (descending) to restore the
#    sheet's original ordering convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the Saldo value for existing accounts -----------------------
# xlWhole (=1) as the LookAt arg forces an exact (non-substring) match on
# the account number.

# 004207955 / SILVANIA : 20       -> 6003.69
$hit = $ws.Range("A1:A184").Find("004207955", $null, $null, 1)
$hit.Offset(0, 2).Value = 6003.69

# 004377713 / DANIELI   : 104.98  -> 496.98
$hit = $ws.Range("A1:A184").Find("004377713", $null, $null, 1)
$hit.Offset(0, 2).Value = 496.98

# 004368628 / CAMILA    : 0.81    -> 363.8
$hit = $ws.Range("A1:A184").Find("004368628", $null, $null, 1)
$hit.Offset(0, 2).Value = 363.8

# --- 2. Insert the new account row -----------------------------------------
# Inserted just above the last data row (row 185 currently holds
# 004976625/NORTON, the final data row before the blank separator/footer).
$newRow = 185
$ws.Rows.Item($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "'004224405"
$ws.Cells.Item($newRow, 2).Value = "MILA"
$ws.Cells.Item($newRow, 3).Value = 5.88

# --- 3. Re-sort the data block (rows 2..185) by Saldo, descending ----------
$dataRange = $ws.Range("A2:C185")
$sortKey = $ws.Range("C2:C185")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)
